$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range('D2').Value = "'30.081.23"
$ws.Range('E2').Value = '  +0.07%  '

# Row 3: Ethereum -> Ethereum
$ws.Range('D3').Value = "'1.912.88"
$ws.Range('E3').Value = '  +0.29%  '

# Row 4: TetherUSD -> TetherUSD
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  +0.09%  '

# Row 5: XRP -> XRP
$ws.Range('D5').Value = "'0.7927"
$ws.Range('E5').Value = '  +5.99%  '

# Row 6: BNB -> BNB
$ws.Range('D6').Value = "'243.12"
$ws.Range('E6').Value = '  +0.21%  '

# Row 8: LidoStakedEther -> Cardano
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.3177"
$ws.Range('E8').Value = '  +2.88%  '

# Row 9: Cardano -> Solana
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').Value = "'26.29"
$ws.Range('E9').Value = '  -0.68%  '

# Row 10: Solana -> Dogecoin
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = "'0.06952"
$ws.Range('E10').Value = '  -0.09%  '

# Row 11: Dogecoin -> TRON
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = "'0.08003"
$ws.Range('E11').Value = '  -1.06%  '

# Row 12: TRON -> Polygon
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = "'0.7505"
$ws.Range('E12').Value = '  -2.23%  '

# Row 13: Polygon -> WrappedEther
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = "'1.914.25"
$ws.Range('E13').Value = '  -0.32%  '

# Row 14: WrappedEther -> Polkadot
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'5.234"
$ws.Range('E14').Value = '  -0.72%  '

# Row 15: Polkadot -> Litecoin
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = "'93.62"
$ws.Range('E15').Value = '  +1.99%  '

# Row 16: Litecoin -> WrappedBTC
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = "'30.093.77"
$ws.Range('E16').Value = '  +0.07%  '

# Row 17: WrappedBTC -> Avalanche
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').Value = "'14.04"
$ws.Range('E17').Value = '  -0.71%  '

# Row 18: Avalanche -> Uniswap
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = "'5.957"
$ws.Range('E18').Value = '  -2.37%  '

# Row 19: Uniswap -> BitcoinCash
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = "'248.10"
$ws.Range('E19').Value = '  +3.49%  '

# Row 20: BitcoinCash -> ShibaInu
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'0.000007809"
$ws.Range('E20').Value = '  +0.02%  '

# Row 21: ShibaInu -> Dai
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = "'1.001"
$ws.Range('E21').Value = '  +0.07%  '

# Row 22: Dai -> BinanceUSD
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = "'1.002"
$ws.Range('E22').Value = '  +0.12%  '

# Row 23: BinanceUSD -> Chainlink
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').Value = "'6.932"
$ws.Range('E23').Value = '  -2.38%  '

# Row 24: Chainlink -> Monero
$ws.Range('B24').Value = 'Monero'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D24').Value = "'168.97"
$ws.Range('E24').Value = '  +1.23%  '

# Row 25: Monero -> Cosmos
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value = "'9.318"
$ws.Range('E25').Value = '  -0.25%  '

# Row 26: Cosmos -> Stellar
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').Value = "'0.1389"
$ws.Range('E26').Value = '  +9.06%  '

# Row 27: Stellar -> EthereumClassic
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = "'18.95"
$ws.Range('E27').Value = '  +0.04%  '

# Row 28: EthereumClassic -> LidoDAOToken
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = "'2.057"
$ws.Range('E28').Value = '  +0.18%  '

# Row 29: LidoDAOToken -> Toncoin
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = "'1.384"
$ws.Range('E29').Value = '  +2.35%  '

# Row 30: Toncoin -> PancakeSwap
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = "'1.524"
$ws.Range('E30').Value = '  -0.61%  '

# Row 31: PancakeSwap -> Filecoin
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'4.348"
$ws.Range('E31').Value = '  +0.70%  '

# Row 32: Filecoin -> Hedera
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.05634"
$ws.Range('E32').Value = '  +4.48%  '

# Row 33: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$ws.Range('D33').Value = "'4.122"
$ws.Range('E33').Value = '  +1.24%  '

# Row 34: Hedera -> ARBITRUM
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').Value = "'1.262"
$ws.Range('E34').Value = '  -3.05%  '

# Row 35: ARBITRUM -> ImmutableX
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = "'0.7377"
$ws.Range('E35').Value = '  -0.67%  '

# Row 36: ImmutableX -> HuobiToken
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = "'2.730"
$ws.Range('E36').Value = '  +0.39%  '

# Row 37: HuobiToken -> VeChain
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.01925"
$ws.Range('E37').Value = '  -2.16%  '

# Row 38: VeChain -> MXToken
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = "'2.796"
$ws.Range('E38').Value = '  -0.11%  '

# Row 39: MXToken -> FraxShare
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = "'6.197"
$ws.Range('E39').Value = '  -1.54%  '

# Row 40: FraxShare -> TheSandbox
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = "'0.4454"
$ws.Range('E40').Value = '  -0.47%  '

# Row 41: TheSandbox -> Aave
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'73.00"
$ws.Range('E41').Value = '  -1.48%  '

# Row 42: Aave -> PaxDollar
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.002"
$ws.Range('E42').Value = '  +0.00%  '

# Row 43: RenderToken -> RenderToken
$ws.Range('D43').Value = "'1.906"
$ws.Range('E43').Value = '  -3.32%  '

# Row 44: PaxDollar -> TrustWalletToken
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = "'0.8329"
$ws.Range('E44').Value = '  -0.35%  '

# Row 45: TrustWalletToken -> Aptos
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = "'7.607"
$ws.Range('E45').Value = '  -1.01%  '

# Row 46: Aptos -> Quant
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').Value = "'101.04"
$ws.Range('E46').Value = '  -0.82%  '

# Row 47: Quant -> EnergySwap
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = "'9.846"
$ws.Range('E47').Value = '  -0.25%  '

# Row 48: EnergySwap -> Maker
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = "'986.92"
$ws.Range('E48').Value = '  +7.19%  '

# Row 49: Maker -> RocketPoolETH
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = "'2.065.87"
$ws.Range('E49').Value = '  +0.35%  '

# Row 50: RocketPoolETH -> Elrond
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = "'36.53"

# Row 51: Elrond -> NEARProtocol
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = "'1.502"
$ws.Range('E51').Value = '  +1.43%  '
